$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the propuesta numbers in column A
# A2 keeps its plain "text" style (no quote prefix), matching the original file.
$ws.Range("A2").Value = "4873593"
# A3 originally carries the "quotePrefix" text style; use a leading apostrophe
# so Excel keeps treating it as text-with-quote-prefix instead of resetting
# the cell style when assigning a numeric-looking string.
$ws.Range("A3").Value = "'4873585"

# Update the selected/active cell to E7 (matches saved view state)
$ws.Range("E7").Select()
